$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Update the "Diferença" period label (2024/03-2023/03 -> 2024/07-2023/07)
#    and refresh the rounded values for rows 2-8 (Roraima .. Sergipe).
# ------------------------------------------------------------------
$newLabel = "Diferença 2024/07 - 2023/07"

$ws.Range("B2").Value = $newLabel
$ws.Range("C2").Value = 3.18

$ws.Range("B3").Value = $newLabel
$ws.Range("C3").Value = 2.76

$ws.Range("B4").Value = $newLabel
$ws.Range("C4").Value = 2.66

$ws.Range("B5").Value = $newLabel
$ws.Range("C5").Value = 2.49

$ws.Range("B6").Value = $newLabel
$ws.Range("C6").Value = 2.32

$ws.Range("B7").Value = $newLabel
$ws.Range("C7").Value = 2.26

$ws.Range("B8").Value = $newLabel
$ws.Range("C8").Value = 1.29

# ------------------------------------------------------------------
# 2. Rows 9/10 swap order (Brasil now comes before Nordeste) and get
#    refreshed values too. Column D on these summary rows is blank.
# ------------------------------------------------------------------
$ws.Range("A9").Value = "Brasil"
$ws.Range("B9").Value = $newLabel
$ws.Range("C9").Value = 0.54
$ws.Range("D9").Value = ""

$ws.Range("A10").Value = "Nordeste"
$ws.Range("B10").Value = $newLabel
$ws.Range("C10").Value = 0.46
$ws.Range("D10").Value = ""

# ------------------------------------------------------------------
# 3. Header row (A1:D1) picks up a thin box border plus top vertical
#    alignment (in addition to the existing bold font + centered text).
# ------------------------------------------------------------------
$headerRange = $ws.Range("A1:D1")
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# ------------------------------------------------------------------
# 4. Page margins reset back to Excel's built-in defaults.
# ------------------------------------------------------------------
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
